# Apply updated TPM-derived values to the LR-pair worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"=35.77148166666667; "H2"=107.314445; "I2"=0.1058641704420874; "J2"=0.1112463097643854; "M2"=0.1332063333333333; "N2"=0.399619; "O2"=0.004085296756603924; "P2"=0.004098425343137321; "Q2"=4.764987910717222; "R2"=42.884891196455; "S2"=0.0004324865521476245; "T2"=0.0004559346952688618
    "G3"=35.77148166666667; "H3"=107.314445; "I3"=0.1058641704420874; "J3"=0.1112463097643854; "O3"=0.0354817324688748; "P3"=0.03559575723202559; "Q3"=41.38500489399888; "R3"=372.46504404599; "S3"=0.003756244173665507; "T3"=0.00395989663533178
    "G4"=35.77148166666667; "H4"=107.314445; "I4"=0.1058641704420874; "J4"=0.1112463097643854; "M4"=16.73711766666667; "N4"=50.211353; "O4"=0.5133096213032781; "P4"=0.5149592027616658; "Q4"=598.7114977660094; "R4"=5388.403479894086; "S4"=0.05434109723921354; "T4"=0.05728731098644522
    "G5"=35.77148166666667; "H5"=107.314445; "I5"=0.1058641704420874; "J5"=0.1112463097643854; "K5"=1; "L5"=0.5; "M5"=0.3133455; "N5"=0.626691; "O5"=0.009609973661260611; "P5"=0.006427237635638122; "Q5"=11.2088328085825; "R5"=67.252996851495; "S5"=0.001017351889619664; "T5"=0.0007150064689435144
    "G6"=35.77148166666667; "H6"=107.314445; "I6"=0.1058641704420874; "J6"=0.1112463097643854; "M6"=14.265684; "N6"=42.797052; "O6"=0.4375133758099827; "P6"=0.4389193770275331; "Q6"=510.30465366846; "R6"=4592.74188301614; "S6"=0.04631699058744103; "T6"=0.048828160978396
    "I7"=0.4187622210170216; "J7"=0.4400521117044616; "M7"=0.1332063333333333; "N7"=0.399619; "O7"=0.004085296756603924; "P7"=0.004098425343137321; "Q7"=18.84865211977245; "R7"=169.637869077952; "S7"=0.001710767943309094; "T7"=0.001803520726910661
    "I8"=0.4187622210170216; "J8"=0.4400521117044616; "O8"=0.0354817324688748; "P8"=0.03559575723202559; "S8"=0.01485840909419778; "T8"=0.01566398813767222
    "I9"=0.4187622210170216; "J9"=0.4400521117044616; "M9"=16.73711766666667; "N9"=50.211353; "O9"=0.5133096213032781; "P9"=0.5149592027616658; "Q9"=2368.296615426425; "R9"=21314.66953883782; "S9"=0.214954677086367; "T9"=0.2266088846169171
    "I10"=0.4187622210170216; "J10"=0.4400521117044616; "K10"=1; "L10"=0.5; "M10"=0.3133455; "N10"=0.626691; "O10"=0.009609973661260611; "P10"=0.006427237635638122; "Q10"=44.338284637088; "R10"=266.029707822528; "S10"=0.004024293914304572; "T10"=0.002828319493988947
    "I11"=0.4187622210170216; "J11"=0.4400521117044616; "M11"=14.265684; "N11"=42.797052; "O11"=0.4375133758099827; "P11"=0.4389193770275331; "Q11"=2018.589568813824; "R11"=18167.30611932442; "S11"=0.1832140729788432; "T11"=0.1931473987289727
    "G12"=52.33127733333333; "H12"=156.993832; "I12"=0.1548721776383825; "J12"=0.1627458863135329; "M12"=0.1332063333333333; "N12"=0.399619; "O12"=0.004085296756603924; "P12"=0.004098425343137321; "Q12"=6.970857572223111; "R12"=62.737718150008; "S12"=0.0006326988049942709; "T12"=0.0006670018649587284
    "G13"=52.33127733333333; "H13"=156.993832; "I13"=0.1548721776383825; "J13"=0.1627458863135329; "O13"=0.0354817324688748; "P13"=0.03559575723202559; "Q13"=60.54348513518044; "R13"=544.891366216624; "S13"=0.005495133173837142; "T13"=0.005793063059727353
    "G14"=52.33127733333333; "H14"=156.993832; "I14"=0.1548721776383825; "J14"=0.1627458863135329; "M14"=16.73711766666667; "N14"=50.211353; "O14"=0.5133096213032781; "P14"=0.5149592027616658; "Q14"=875.8747463749662; "R14"=7882.872717374697; "S14"=0.07949737885397214; "T14"=0.08380749186875759
    "G15"=52.33127733333333; "H15"=156.993832; "I15"=0.1548721776383825; "J15"=0.1627458863135329; "K15"=1; "L15"=0.5; "M15"=0.3133455; "N15"=0.626691; "O15"=0.009609973661260611; "P15"=0.006427237635638122; "Q15"=16.397770261652; "R15"=98.38662156991199; "S15"=0.001488317547966931; "T15"=0.001046006485559622
    "G16"=52.33127733333333; "H16"=156.993832; "I16"=0.1548721776383825; "J16"=0.1627458863135329; "M16"=14.265684; "N16"=42.797052; "O16"=0.4375133758099827; "P16"=0.4389193770275331; "Q16"=746.541465753696; "R16"=6718.873191783264; "S16"=0.06775864925761205; "T16"=0.07143232303452957
    "G17"=49.043167; "H17"=98.08633399999999; "I17"=0.1451411557029742; "J17"=0.1016800925151965; "M17"=0.1332063333333333; "N17"=0.399619; "O17"=0.004085296756603924; "P17"=0.004098425343137321; "Q17"=6.532860451124334; "R17"=39.197162706746; "S17"=0.0005929446926431056; "T17"=0.0004167282680568286
    "G18"=49.043167; "H18"=98.08633399999999; "I18"=0.1451411557029742; "J18"=0.1016800925151965; "O18"=0.0354817324688748; "P18"=0.03559575723202559; "Q18"=56.73938041553132; "R18"=340.4362824931879; "S18"=0.005149859656876232; "T18"=0.003619379888500836
    "G19"=49.043167; "H19"=98.08633399999999; "I19"=0.1451411557029742; "J19"=0.1016800925151965; "M19"=16.73711766666667; "N19"=50.211353; "O19"=0.5133096213032781; "P19"=0.5149592027616658; "Q19"=820.8412568249836; "R19"=4925.047540949902; "S19"=0.07450235166941381; "T19"=0.052361099378358
    "G20"=49.043167; "H20"=98.08633399999999; "I20"=0.1451411557029742; "J20"=0.1016800925151965; "K20"=1; "L20"=0.5; "M20"=0.3133455; "N20"=0.626691; "O20"=0.009609973661260611; "P20"=0.006427237635638122; "Q20"=15.3674556851985; "R20"=61.46982274079399; "S20"=0.001394802683470508; "T20"=0.0006535221174088368
    "G21"=49.043167; "H21"=98.08633399999999; "I21"=0.1451411557029742; "J21"=0.1016800925151965; "M21"=14.265684; "N21"=42.797052; "O21"=0.4375133758099827; "P21"=0.4389193770275331; "Q21"=699.634322781228; "R21"=4197.805936687368; "S21"=0.06350119700057057; "T21"=0.04462936286287197
    "G22"=59.25420133333333; "H22"=177.762604; "I22"=0.1753602751995342; "J22"=0.1842755997024237; "M22"=0.1332063333333333; "N22"=0.399619; "O22"=0.004085296756603924; "P22"=0.004098425343137321; "Q22"=7.893034894208445; "R22"=71.03731404787601; "S22"=0.0007163987635098287; "T22"=0.0007552397879422416
    "G23"=59.25420133333333; "H23"=177.762604; "I23"=0.1753602751995342; "J23"=0.1842755997024237; "O23"=0.0354817324688748; "P23"=0.03559575723202559; "Q23"=68.55280513736977; "R23"=616.975246236328; "S23"=0.006222086370298133; "T23"=0.006559429510793403
    "G24"=59.25420133333333; "H24"=177.762604; "I24"=0.1753602751995342; "J24"=0.1842755997024237; "M24"=16.73711766666667; "N24"=50.211353; "O24"=0.5133096213032781; "P24"=0.5149592027616658; "Q24"=991.7445399603569; "R24"=8925.700859643213; "S24"=0.09001411645431154; "T24"=0.094894415911188
    "G25"=59.25420133333333; "H25"=177.762604; "I25"=0.1753602751995342; "J25"=0.1842755997024237; "K25"=1; "L25"=0.5; "M25"=0.3133455; "N25"=0.626691; "O25"=0.009609973661260611; "P25"=0.006427237635638122; "Q25"=18.567037343894; "R25"=111.402224063364; "S25"=0.001685207625898936; "T25"=0.001184383069737203
    "G26"=59.25420133333333; "H26"=177.762604; "I26"=0.1753602751995342; "J26"=0.1842755997024237; "M26"=14.265684; "N26"=42.797052; "O26"=0.4375133758099827; "P26"=0.4389193770275331; "Q26"=845.3017118937121; "R26"=7607.715407043409; "S26"=0.07672246598551581; "T26"=0.08088213142276289
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
